$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.176.36'
$ws.Range("E2").Value = '  -4.26%  '
$ws.Range("D3").Value = '2.972.86'
$ws.Range("E3").Value = '  -3.84%  '
$ws.Range("E4").Value = '  -0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '543.78'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '130.69'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -6.93%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '2.970.27'
$ws.Range("E8").Value = '  -3.74%  '
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("E10").Value = '  -7.44%  '
$ws.Range("E11").Value = '  -9.76%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.441'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.55%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000218'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.65%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '33.47'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.23%  '
$ws.Range("D15").Value = '3.449.94'
$ws.Range("E15").Value = '  -4.02%  '
$ws.Range("D16").Value = '61.168.93'
$ws.Range("E16").Value = '  -4.40%  '
$ws.Range("E17").Value = '  -3.06%  '
$ws.Range("D18").Value = '2.970.08'
$ws.Range("E18").Value = '  -3.99%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.53'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -2.09%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '465.45'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.98%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.93'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.99%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.660'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -5.87%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.90'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.09%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '79.33'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("E25").Value = '  -4.17%  '
$ws.Range("E26").Value = '  +0.06%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.68'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.63%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.54'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -6.31%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.07%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '25.23'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -4.31%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.12'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.20%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.26'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.08%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.39'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.89%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '54.44'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -4.87%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.79'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -3.92%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '443.76'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -10.85%  '
$ws.Range("D38").Value = '3.129.52'
$ws.Range("E38").Value = '  -3.55%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0781'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0374'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -7.39%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.02%  '
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("E43").Value = '  -0.03%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -12.57%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '25.31'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.01%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.238'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -6.26%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.107'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.92'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -5.67%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '114.81'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -7.76%  '
$ws.Range("E50").Value = '  +7.99%  '
$ws.Range("D51").Value = '0.0₃0475'
$ws.Range("E51").Value = '  -11.06%  '
